$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.945.41"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.828.75"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.11"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.320"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.06%  "
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0991"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "2.091.11"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").Value = "1.822.74"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.669"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "34.951.26"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "0.0₃0786"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  -5.71%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.701"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "92.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "1.338.78"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0193"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").Value = "2.006.70"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0669"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.13%  "
